$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 3407.3333
$ws.Range("I39").Value = 194.33333
$ws.Range("J39").Value = 9833.333000000001
$ws.Range("K39").Value = 582.99999
$ws.Range("L39").Value = 29499.999
$ws.Range("M39").Value = -286.99999
$ws.Range("N39").Value = -30091.999

$ws.Range("H64").Value = 50002748
$ws.Range("I64").Value = 66668664
$ws.Range("J64").Value = 4999
$ws.Range("K64").Value = 66668664
$ws.Range("L64").Value = 4999
$ws.Range("M64").Value = -66668416
$ws.Range("N64").Value = -5495

$ws.Range("H67").Value = 50002748
$ws.Range("I67").Value = 66668664
$ws.Range("J67").Value = 4999
$ws.Range("K67").Value = 66668664
$ws.Range("L67").Value = 4999
$ws.Range("M67").Value = -66667806
$ws.Range("N67").Value = -6715

$ws.Range("H69").Value = 15233.294
$ws.Range("J69").Value = 15931.934
$ws.Range("L69").Value = 47795.802
$ws.Range("N69").Value = -49543.802

$ws.Range("H72").Value = 15233.294
$ws.Range("J72").Value = 15931.934
$ws.Range("L72").Value = 143387.406
$ws.Range("N72").Value = -152123.406

$ws.Range("H76").Value = 4798.2
$ws.Range("I76").Value = 4798.2
$ws.Range("K76").Value = 4798.2
$ws.Range("M76").Value = -4483.2

$ws.Range("H79").Value = 4798.2
$ws.Range("I79").Value = 4798.2
$ws.Range("K79").Value = 4798.2
$ws.Range("M79").Value = -3706.2

$ws.Range("H113").Value = 4402
$ws.Range("I113").Value = 3600
$ws.Range("K113").Value = 3600
$ws.Range("M113").Value = -346

$ws.Range("H132").Value = 2033.2572
$ws.Range("I132").Value = 2060.7058
$ws.Range("K132").Value = 6182.117400000001
$ws.Range("M132").Value = -3652.117400000001

$ws.Range("H137").Value = 1936.119
$ws.Range("I137").Value = 1325.4642
$ws.Range("J137").Value = 3157.4285
$ws.Range("K137").Value = 3976.3926
$ws.Range("L137").Value = 9472.2855
$ws.Range("M137").Value = -1426.3926
$ws.Range("N137").Value = -14572.2855

$ws.Range("H138").Value = 3631.7126
$ws.Range("J138").Value = 3649.419
$ws.Range("L138").Value = 10948.257
$ws.Range("N138").Value = -21228.257

$ws.Range("H141").Value = 3547.087
$ws.Range("I141").Value = 2052.1667
$ws.Range("J141").Value = 5177.909
$ws.Range("K141").Value = 6156.500100000001
$ws.Range("L141").Value = 15533.727
$ws.Range("M141").Value = -976.5001000000011
$ws.Range("N141").Value = -25893.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10003.087
$ws.Range("I32").Value = 5363.476
$ws.Range("J32").Value = 17220.26
$ws.Range("K32").Value = 5363.476
$ws.Range("L32").Value = 17220.26
$ws.Range("M32").Value = -5076.476
$ws.Range("N32").Value = -17794.26

$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9533

$ws.Range("H45").Value = 1055.2858
$ws.Range("I45").Value = 932.0833
$ws.Range("K45").Value = 932.0833
$ws.Range("M45").Value = -555.0833

$ws.Range("H61").Value = 33336164
$ws.Range("I61").Value = 35716908
$ws.Range("J61").Value = 5740
$ws.Range("K61").Value = 35716908
$ws.Range("L61").Value = 5740
$ws.Range("M61").Value = -35716696
$ws.Range("N61").Value = -6164

$ws.Range("H97").Value = 540.5
$ws.Range("I97").Value = 604
$ws.Range("J97").Value = 159.5
$ws.Range("K97").Value = 604
$ws.Range("L97").Value = 159.5
$ws.Range("M97").Value = -108
$ws.Range("N97").Value = -1151.5

$ws.Range("H132").Value = 2506980.5
$ws.Range("I132").Value = 3131644.5
$ws.Range("K132").Value = 9394933.5
$ws.Range("M132").Value = -9392403.5

$ws.Range("H136").Value = 33336164
$ws.Range("I136").Value = 35716908
$ws.Range("J136").Value = 5740
$ws.Range("K136").Value = 107150724
$ws.Range("L136").Value = 17220
$ws.Range("M136").Value = -107148174
$ws.Range("N136").Value = -22320

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2202.9167
$ws.Range("I20").Value = 2173.4285
$ws.Range("J20").Value = 2244.2
$ws.Range("K20").Value = 2173.4285
$ws.Range("L20").Value = 2244.2
$ws.Range("M20").Value = -1926.4285
$ws.Range("N20").Value = -2738.2

$ws.Range("H22").Value = 47619050
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 47619050
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 47619050
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -47619396

$ws.Range("H94").Value = 700.3333
$ws.Range("I94").Value = 687.3
$ws.Range("J94").Value = 726.4
$ws.Range("K94").Value = 687.3
$ws.Range("L94").Value = 726.4
$ws.Range("M94").Value = -236.3
$ws.Range("N94").Value = -1628.4

$ws.Range("H134").Value = 13160436
$ws.Range("J134").Value = 4981.3335
$ws.Range("L134").Value = 14944.0005
$ws.Range("N134").Value = -20014.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7753.0713
$ws.Range("I22").Value = 9549.454
$ws.Range("J22").Value = 1166.3334
$ws.Range("K22").Value = 9549.454
$ws.Range("L22").Value = 1166.3334
$ws.Range("M22").Value = -9199.454
$ws.Range("N22").Value = -1866.3334

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H99").Value = 5741.2856
$ws.Range("I99").Value = 6599.75
$ws.Range("K99").Value = 6599.75
$ws.Range("M99").Value = -5101.75

$ws.Range("H126").Value = 5741.2856
$ws.Range("I126").Value = 6599.75
$ws.Range("K126").Value = 19799.25
$ws.Range("M126").Value = -17329.25

$ws.Range("H127").Value = 104995
$ws.Range("J127").Value = 104995
$ws.Range("L127").Value = 104995
$ws.Range("N127").Value = -114915

$ws.Range("H134").Value = 41670936
$ws.Range("I134").Value = 50003124
$ws.Range("K134").Value = 150009372
$ws.Range("M134").Value = -150006837

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 160.88889
$ws.Range("J17").Value = 149.75
$ws.Range("L17").Value = 449.25
$ws.Range("N17").Value = -787.25

$ws.Range("H34").Value = 840.0909
$ws.Range("I34").Value = 834.1
$ws.Range("J34").Value = 900
$ws.Range("K34").Value = 2502.3
$ws.Range("L34").Value = 2700
$ws.Range("M34").Value = -2418.3
$ws.Range("N34").Value = -2868

$ws.Range("H39").Value = 925
$ws.Range("I39").Value = 925
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2775
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2481
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 2911.111
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 2900
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 8700
$ws.Range("M55").Value = -8823
$ws.Range("N55").Value = -9054

$ws.Range("H92").Value = 337.4
$ws.Range("I92").Value = 337.4
$ws.Range("K92").Value = 1012.2
$ws.Range("M92").Value = 235.8000000000001

$ws.Range("H108").Value = 1361.2727
$ws.Range("I108").Value = 1197.5
$ws.Range("K108").Value = 3592.5
$ws.Range("M108").Value = -712.5

$ws.Range("H113").Value = 51240.95
$ws.Range("I113").Value = 126230.875
$ws.Range("J113").Value = 1247.6666
$ws.Range("K113").Value = 378692.625
$ws.Range("L113").Value = 3742.9998
$ws.Range("M113").Value = -376522.625
$ws.Range("N113").Value = -8082.9998

$ws.Range("H122").Value = 1587.4117
$ws.Range("J122").Value = 2770.1667
$ws.Range("L122").Value = 24931.5003
$ws.Range("N122").Value = -29831.5003

$ws.Range("H131").Value = 1827.9429
$ws.Range("J131").Value = 1799.3334
$ws.Range("L131").Value = 5398.0002
$ws.Range("N131").Value = -15478.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 886.6923
$ws.Range("J97").Value = 1121.6666
$ws.Range("L97").Value = 1121.6666
$ws.Range("N97").Value = -2113.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5472.826
$ws.Range("I61").Value = 5108.8
$ws.Range("K61").Value = 5108.8
$ws.Range("M61").Value = -4906.8

$ws.Range("H113").Value = 5472.826
$ws.Range("I113").Value = 5108.8
$ws.Range("K113").Value = 5108.8
$ws.Range("M113").Value = -2938.8

$ws.Range("H132").Value = 88768290
$ws.Range("I132").Value = 101445040
$ws.Range("K132").Value = 304335120
$ws.Range("M132").Value = -304332590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11158
$ws.Range("I41").Value = 11499
$ws.Range("J41").Value = 10987.5
$ws.Range("K41").Value = 11499
$ws.Range("L41").Value = 10987.5
$ws.Range("M41").Value = -11109
$ws.Range("N41").Value = -11767.5

$ws.Range("H132").Value = 29422838
$ws.Range("I132").Value = 50007028
$ws.Range("J132").Value = 16855.857
$ws.Range("K132").Value = 150021084
$ws.Range("L132").Value = 50567.571
$ws.Range("M132").Value = -150018554
$ws.Range("N132").Value = -55627.571
